$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08943201502788156
$ws.Range("D2").Value = 0.1948903110990585
$ws.Range("E2").Value = 0.1678563450209651
$ws.Range("F2").Value = 1.381042796064705
$ws.Range("G2").Value = 0.002451400755762983
$ws.Range("I2").Value = 0.4533284278391356
$ws.Range("J2").Value = 0.186868640309811
$ws.Range("K2").Value = 0.8053691607786959
$ws.Range("N2").Value = 1.371772158472872
$ws.Range("O2").Value = 3.284725044647814
$ws.Range("B3").Value = 0.07898517216182199
$ws.Range("D3").Value = 0.1898770388727229
$ws.Range("E3").Value = 0.1638719165803231
$ws.Range("F3").Value = 1.377974351892817
$ws.Range("G3").Value = 0.002454188753360244
$ws.Range("I3").Value = 0.4596567733388368
$ws.Range("J3").Value = 0.1826903552363603
$ws.Range("K3").Value = 0.7210469065795735
$ws.Range("N3").Value = 1.384201999929516
$ws.Range("O3").Value = 3.288251460566727
$ws.Range("B4").Value = 0.07256721165369129
$ws.Range("D4").Value = 0.1868784677719333
$ws.Range("E4").Value = 0.1615098194224771
$ws.Range("F4").Value = 1.376867550056232
$ws.Range("G4").Value = 0.002455992836384807
$ws.Range("I4").Value = 0.4637650084354483
$ws.Range("J4").Value = 0.180231683285065
$ws.Range("K4").Value = 0.6693078323392854
$ws.Range("N4").Value = 1.39231691949788
$ws.Range("O4").Value = 3.292336158825208
$ws.Range("B5").Value = 0.06995116669473589
$ws.Range("D5").Value = 0.1856766324671355
$ws.Range("E5").Value = 0.1605684965687466
$ws.Range("F5").Value = 1.376611969870353
$ws.Range("G5").Value = 0.002456751279359963
$ws.Range("I5").Value = 0.465495125160539
$ws.Range("J5").Value = 0.1792566179308821
$ws.Range("K5").Value = 0.6482336586487349
$ws.Range("N5").Value = 1.395745341284627
$ws.Range("O5").Value = 3.294483170535614
$ws.Range("B6").Value = 0.06951674020773169
$ws.Range("D6").Value = 0.1854782857132449
$ws.Range("E6").Value = 0.160413475263951
$ws.Range("F6").Value = 1.376581335997194
$ws.Range("G6").Value = 0.002456878625644036
$ws.Range("I6").Value = 0.4657857906824976
$ws.Range("J6").Value = 0.1790963317034837
$ws.Range("K6").Value = 0.6447349392554713
$ws.Range("N6").Value = 1.396321970554201
$ws.Range("O6").Value = 3.294868814036789
$ws.Range("B7").Value = 0.07253193319225204
$ws.Range("D7").Value = 0.1868621778680506
$ws.Range("E7").Value = 0.1614970383176946
$ws.Range("F7").Value = 1.376863311831244
$ws.Range("G7").Value = 0.002456002970823023
$ws.Range("I7").Value = 0.4637881147411056
$ws.Range("J7").Value = 0.18021842443693
$ws.Range("K7").Value = 0.6690235770620347
$ws.Range("N7").Value = 1.392362664224265
$ws.Range("O7").Value = 3.292363160969984
$ws.Range("B8").Value = 0.08583081072562493
$ws.Range("D8").Value = 0.1931452669322624
$ws.Range("E8").Value = 0.1664650270109362
$ws.Range("F8").Value = 1.37982348125901
$ws.Range("G8").Value = 0.002452342954775938
$ws.Range("I8").Value = 0.4554642412873857
$ws.Range("J8").Value = 0.1854057998214103
$ws.Range("K8").Value = 0.7762881807904591
$ws.Range("N8").Value = 1.375957783609682
$ws.Range("O8").Value = 3.285542437972993
$ws.Range("B9").Value = 0.1118732146061205
$ws.Range("D9").Value = 0.2060947520767797
$ws.Range("E9").Value = 0.1768755827011432
$ws.Range("F9").Value = 1.391797336870354
$ws.Range("G9").Value = 0.002445894385233122
$ws.Range("I9").Value = 0.4409071224823267
$ws.Range("J9").Value = 0.1964263054521638
$ws.Range("K9").Value = 0.9868751866924583
$ws.Range("N9").Value = 1.347615867606414
$ws.Range("O9").Value = 3.287410407507537
$ws.Range("B10").Value = 0.1309747286336034
$ws.Range("D10").Value = 0.2159886638439872
$ws.Range("E10").Value = 0.1849313376655815
$ws.Range("F10").Value = 1.404361389032047
$ws.Range("G10").Value = 0.002441596403276099
$ws.Range("I10").Value = 0.4312878458171641
$ws.Range("J10").Value = 0.2050420075412802
$ws.Range("K10").Value = 1.141705824701091
$ws.Range("N10").Value = 1.329121076966686
$ws.Range("O10").Value = 3.298100095617883
$ws.Range("B11").Value = 0.1396556875444617
$ws.Range("D11").Value = 0.2205715622228723
$ws.Range("E11").Value = 0.1886845061943916
$ws.Range("F11").Value = 1.410896580350737
$ws.Range("G11").Value = 0.002439735686347182
$ws.Range("I11").Value = 0.4271452282810992
$ws.Range("J11").Value = 0.2090746770264929
$ws.Range("K11").Value = 1.212160241518689
$ws.Range("N11").Value = 1.321211679004819
$ws.Range("O11").Value = 3.304991988469538
$ws.Range("B12").Value = 0.1429415300560066
$ws.Range("D12").Value = 0.2223187176825121
$ws.Range("E12").Value = 0.1901184456545835
$ws.Range("F12").Value = 1.413489226959115
$ws.Range("G12").Value = 0.002439044592226013
$ws.Range("I12").Value = 0.4256100650264969
$ws.Range("J12").Value = 0.2106180570014118
$ws.Range("K12").Value = 1.238841648763696
$ws.Range("N12").Value = 1.318288997887521
$ws.Range("O12").Value = 3.307893909993197
$ws.Range("B13").Value = 0.1422339332759321
$ws.Range("D13").Value = 0.2219419171554193
$ws.Range("E13").Value = 0.1898090571097271
$ws.Range("F13").Value = 1.412925609440336
$ws.Range("G13").Value = 0.002439192831447719
$ws.Range("I13").Value = 0.4259391972199484
$ws.Range("J13").Value = 0.2102849378114087
$ws.Range("K13").Value = 1.233095267151953
$ws.Range("N13").Value = 1.318915228457506
$ws.Range("O13").Value = 3.307255932293401
$ws.Range("B14").Value = 0.1399260461722633
$ws.Range("D14").Value = 0.2207150675899072
$ws.Range("E14").Value = 0.1888022230374062
$ws.Range("F14").Value = 1.411107515632438
$ws.Range("G14").Value = 0.002439678559052655
$ws.Range("I14").Value = 0.4270182568667771
$ws.Range("J14").Value = 0.209201325256771
$ws.Range("K14").Value = 1.214355306754612
$ws.Range("N14").Value = 1.320969776727843
$ws.Range("O14").Value = 3.305224875067324
$ws.Range("B15").Value = 0.1385122030220742
$ws.Range("D15").Value = 0.2199651091213894
$ws.Range("E15").Value = 0.1881871603122249
$ws.Range("F15").Value = 1.410009237236579
$ws.Range("G15").Value = 0.002439977839821776
$ws.Range("I15").Value = 0.4276835825453658
$ws.Range("J15").Value = 0.2085397033468723
$ws.Range("K15").Value = 1.202876745141566
$ws.Range("N15").Value = 1.322237679227321
$ws.Range("O15").Value = 3.304018845160527
$ws.Range("B16").Value = 0.1304072178486706
$ws.Range("D16").Value = 0.215690804700742
$ws.Range("E16").Value = 0.1846878380545931
$ws.Range("F16").Value = 1.403950797076845
$ws.Range("G16").Value = 0.002441719900906392
$ws.Range("I16").Value = 0.4315632704027053
$ws.Range("J16").Value = 0.2047807436070457
$ws.Range("K16").Value = 1.137101786731307
$ws.Range("N16").Value = 1.32964811294471
$ws.Range("O16").Value = 3.297690553014746
$ws.Range("B17").Value = 0.1254327469898868
$ws.Range("D17").Value = 0.2130896191674481
$ws.Range("E17").Value = 0.1825637727568861
$ws.Range("F17").Value = 1.40044411959623
$ws.Range("G17").Value = 0.002442812746099368
$ws.Range("I17").Value = 0.4340030832179638
$ws.Range("J17").Value = 0.202503771398483
$ws.Range("K17").Value = 1.096755553217463
$ws.Range("N17").Value = 1.334323228336821
$ws.Range("O17").Value = 3.294328285971602
$ws.Range("B18").Value = 0.1225707850047115
$ws.Range("D18").Value = 0.2116012201992845
$ws.Range("E18").Value = 0.1813504051300612
$ws.Range("F18").Value = 1.398504327330343
$ws.Range("G18").Value = 0.002443450215865
$ws.Range("I18").Value = 0.4354283516424076
$ws.Range("J18").Value = 0.2012047858253965
$ws.Range("K18").Value = 1.073551532145927
$ws.Range("N18").Value = 1.337059669440805
$ws.Range("O18").Value = 3.292585367656358
$ws.Range("B19").Value = 0.1216016475052015
$ws.Range("D19").Value = 0.2110986051467165
$ws.Range("E19").Value = 0.1809410130574207
$ws.Range("F19").Value = 1.397860797571383
$ws.Range("G19").Value = 0.002443667581945285
$ws.Range("I19").Value = 0.4359146929251967
$ws.Range("J19").Value = 0.2007668040648412
$ws.Range("K19").Value = 1.065695440545113
$ws.Range("N19").Value = 1.33799433008307
$ws.Range("O19").Value = 3.292028034524208
$ws.Range("B20").Value = 0.1259623701418917
$ws.Range("D20").Value = 0.2133657200463261
$ws.Range("E20").Value = 0.1827890204116329
$ws.Range("F20").Value = 1.400809425810621
$ws.Range("G20").Value = 0.002442695490573617
$ws.Range("I20").Value = 0.4337410888453048
$ws.Range("J20").Value = 0.2027450548379335
$ws.Range("K20").Value = 1.101050273575709
$ws.Range("N20").Value = 1.333820644887986
$ws.Range("O20").Value = 3.294666438098119
$ws.Range("B21").Value = 0.1406039697581889
$ws.Range("D21").Value = 0.2210751059693479
$ws.Range("E21").Value = 0.1890976103611663
$ws.Range("F21").Value = 1.411638333502097
$ws.Range("G21").Value = 0.002439535522693922
$ws.Range("I21").Value = 0.4267004000564043
$ws.Range("J21").Value = 0.2095191662997848
$ws.Range("K21").Value = 1.219859646750876
$ws.Range("N21").Value = 1.320364340460152
$ws.Range("O21").Value = 3.305813515365713
$ws.Range("B22").Value = 0.1501645499654103
$ws.Range("D22").Value = 0.2261818533023501
$ws.Range("E22").Value = 0.193294622061515
$ws.Range("F22").Value = 1.4194029447509
$ws.Range("G22").Value = 0.002437549071148534
$ws.Range("I22").Value = 0.4222945013282091
$ws.Range("J22").Value = 0.2140414211587682
$ws.Range("K22").Value = 1.297518660991784
$ws.Range("N22").Value = 1.311992067604841
$ws.Range("O22").Value = 3.314801576033204
$ws.Range("B23").Value = 0.1450627493927215
$ws.Range("D23").Value = 0.2234500758515736
$ws.Range("E23").Value = 0.1910478415527948
$ws.Range("F23").Value = 1.415195926990776
$ws.Range("G23").Value = 0.002438602091150168
$ws.Range("I23").Value = 0.4246281096333639
$ws.Range("J23").Value = 0.2116191188985255
$ws.Range("K23").Value = 1.25607004361882
$ws.Range("N23").Value = 1.316421886516075
$ws.Range("O23").Value = 3.309848563552208
$ws.Range("B24").Value = 0.1257229339519341
$ws.Range("D24").Value = 0.2132408728502071
$ws.Range("E24").Value = 0.1826871616929111
$ws.Range("F24").Value = 1.400644033378001
$ws.Range("G24").Value = 0.002442748473177736
$ws.Range("I24").Value = 0.4338594661441375
$ws.Range("J24").Value = 0.202635939210424
$ws.Range("K24").Value = 1.099108656681892
$ws.Range("N24").Value = 1.334047711303384
$ws.Range("O24").Value = 3.294512967406206
$ws.Range("B25").Value = 0.1048329617678831
$ws.Range("D25").Value = 0.2025246271421963
$ws.Range("E25").Value = 0.1739877333445179
$ws.Range("F25").Value = 1.387897032500604
$ws.Range("G25").Value = 0.002447561343438931
$ws.Range("I25").Value = 0.4446562272943479
$ws.Range("J25").Value = 0.1933539743873785
$ws.Range("K25").Value = 0.9298835799406504
$ws.Range("N25").Value = 1.354873888667662
$ws.Range("O25").Value = 3.285270634156518

Write-Output "Applied 240 cell updates"